$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns at D:E, shifting existing D:K data to F:M
$ws.Columns("D:E").Insert(1)

# Copy number formats from F/G (old D/E) into new D/E so styles match
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(7, 4).Value2 = 43465
$ws.Cells.Item(7, 5).Value2 = 43373
$ws.Cells.Item(8, 4).Value2 = 858200
$ws.Cells.Item(8, 5).Value2 = 843200
$ws.Cells.Item(9, 4).Value2 = 594900
$ws.Cells.Item(9, 5).Value2 = 580700
$ws.Cells.Item(10, 4).Value2 = 263300
$ws.Cells.Item(10, 5).Value2 = 262500
$ws.Cells.Item(12, 4).Value2 = "NA"
$ws.Cells.Item(12, 5).Value2 = "NA"
$ws.Cells.Item(13, 4).Value2 = 0
$ws.Cells.Item(13, 5).Value2 = 0
$ws.Cells.Item(14, 4).Value2 = 0
$ws.Cells.Item(14, 5).Value2 = 2500
$ws.Cells.Item(15, 4).Value2 = 77900
$ws.Cells.Item(15, 5).Value2 = 73100
$ws.Cells.Item(17, 4).Value2 = 816700
$ws.Cells.Item(17, 5).Value2 = 779900
$ws.Cells.Item(18, 4).Value2 = 41500
$ws.Cells.Item(18, 5).Value2 = 63300
$ws.Cells.Item(20, 4).Value2 = -3200
$ws.Cells.Item(20, 5).Value2 = -300
$ws.Cells.Item(21, 4).Value2 = 116200
$ws.Cells.Item(21, 5).Value2 = 136100
$ws.Cells.Item(22, 4).Value2 = 21000
$ws.Cells.Item(22, 5).Value2 = 20700
$ws.Cells.Item(23, 4).Value2 = 17300
$ws.Cells.Item(23, 5).Value2 = 42400
$ws.Cells.Item(24, 4).Value2 = 1100
$ws.Cells.Item(24, 5).Value2 = 11300
$ws.Cells.Item(25, 4).Value2 = 0
$ws.Cells.Item(25, 5).Value2 = 0
$ws.Cells.Item(26, 4).Value2 = 16100
$ws.Cells.Item(26, 5).Value2 = 31100
$ws.Cells.Item(27, 4).Value2 = 16100
$ws.Cells.Item(27, 5).Value2 = 31100
$ws.Cells.Item(28, 4).Value2 = 0
$ws.Cells.Item(28, 5).Value2 = 0
$ws.Cells.Item(29, 4).Value2 = 300
$ws.Cells.Item(29, 5).Value2 = "NA"
$ws.Cells.Item(30, 4).Value2 = 0
$ws.Cells.Item(30, 5).Value2 = 0
$ws.Cells.Item(31, 4).Value2 = 0
$ws.Cells.Item(31, 5).Value2 = 0
$ws.Cells.Item(32, 4).Value2 = 3200
$ws.Cells.Item(32, 5).Value2 = 300
$ws.Cells.Item(33, 4).Value2 = 16400
$ws.Cells.Item(33, 5).Value2 = 31100
$ws.Cells.Item(34, 4).Value2 = 0
$ws.Cells.Item(34, 5).Value2 = 0
$ws.Cells.Item(35, 4).Value2 = 16400
$ws.Cells.Item(35, 5).Value2 = 31100
$ws.Cells.Item(38, 4).Value2 = 43465
$ws.Cells.Item(38, 5).Value2 = 43373
$ws.Cells.Item(41, 4).Value2 = 226500
$ws.Cells.Item(41, 5).Value2 = 215500
$ws.Cells.Item(42, 4).Value2 = 52900
$ws.Cells.Item(42, 5).Value2 = 37400
$ws.Cells.Item(43, 4).Value2 = 661700
$ws.Cells.Item(43, 5).Value2 = 672600
$ws.Cells.Item(44, 4).Value2 = 199500
$ws.Cells.Item(44, 5).Value2 = 196000
$ws.Cells.Item(45, 4).Value2 = 61600
$ws.Cells.Item(45, 5).Value2 = 55300
$ws.Cells.Item(46, 4).Value2 = 1202200
$ws.Cells.Item(46, 5).Value2 = 1176800
$ws.Cells.Item(47, 4).Value2 = 0
$ws.Cells.Item(47, 5).Value2 = 0
$ws.Cells.Item(48, 4).Value2 = 1562000
$ws.Cells.Item(48, 5).Value2 = 1614400
$ws.Cells.Item(49, 4).Value2 = 956100
$ws.Cells.Item(49, 5).Value2 = 965500
$ws.Cells.Item(50, 4).Value2 = 0
$ws.Cells.Item(50, 5).Value2 = 0
$ws.Cells.Item(51, 4).Value2 = 0
$ws.Cells.Item(51, 5).Value2 = 0
$ws.Cells.Item(52, 4).Value2 = 18100
$ws.Cells.Item(52, 5).Value2 = 17600
$ws.Cells.Item(53, 4).Value2 = 0
$ws.Cells.Item(53, 5).Value2 = 0
$ws.Cells.Item(54, 4).Value2 = 3738300
$ws.Cells.Item(54, 5).Value2 = 3774300
$ws.Cells.Item(57, 4).Value2 = 276500
$ws.Cells.Item(57, 5).Value2 = 248400
$ws.Cells.Item(58, 4).Value2 = 7500
$ws.Cells.Item(58, 5).Value2 = 7500
$ws.Cells.Item(59, 4).Value2 = 318300
$ws.Cells.Item(59, 5).Value2 = 320400
$ws.Cells.Item(60, 4).Value2 = 602300
$ws.Cells.Item(60, 5).Value2 = 576300
$ws.Cells.Item(61, 4).Value2 = 1565000
$ws.Cells.Item(61, 5).Value2 = 1616200
$ws.Cells.Item(62, 4).Value2 = 401300
$ws.Cells.Item(62, 5).Value2 = 384500
$ws.Cells.Item(63, 4).Value2 = 0
$ws.Cells.Item(63, 5).Value2 = 0
$ws.Cells.Item(64, 4).Value2 = 0
$ws.Cells.Item(64, 5).Value2 = 0
$ws.Cells.Item(65, 4).Value2 = 0
$ws.Cells.Item(65, 5).Value2 = 0
$ws.Cells.Item(66, 4).Value2 = 2568600
$ws.Cells.Item(66, 5).Value2 = 2577000
$ws.Cells.Item(68, 4).Value2 = 0
$ws.Cells.Item(68, 5).Value2 = 0
$ws.Cells.Item(69, 4).Value2 = 0
$ws.Cells.Item(69, 5).Value2 = 0
$ws.Cells.Item(70, 4).Value2 = 0
$ws.Cells.Item(70, 5).Value2 = 0
$ws.Cells.Item(71, 4).Value2 = 0
$ws.Cells.Item(71, 5).Value2 = 0
$ws.Cells.Item(72, 4).Value2 = 737200
$ws.Cells.Item(72, 5).Value2 = 720700
$ws.Cells.Item(73, 4).Value2 = 0
$ws.Cells.Item(73, 5).Value2 = 0
$ws.Cells.Item(74, 4).Value2 = 0
$ws.Cells.Item(74, 5).Value2 = 0
$ws.Cells.Item(75, 4).Value2 = 0
$ws.Cells.Item(75, 5).Value2 = 0
$ws.Cells.Item(76, 4).Value2 = 1169800
$ws.Cells.Item(76, 5).Value2 = 1197300
$ws.Cells.Item(77, 4).Value2 = 0
$ws.Cells.Item(77, 5).Value2 = 0
$ws.Cells.Item(80, 4).Value2 = 43465
$ws.Cells.Item(80, 5).Value2 = 43373
$ws.Cells.Item(81, 4).Value2 = 16400
$ws.Cells.Item(81, 5).Value2 = 31100
$ws.Cells.Item(83, 4).Value2 = 77900
$ws.Cells.Item(83, 5).Value2 = 73100
$ws.Cells.Item(84, 4).Value2 = 0
$ws.Cells.Item(84, 5).Value2 = 0
$ws.Cells.Item(85, 4).Value2 = 0
$ws.Cells.Item(85, 5).Value2 = 0
$ws.Cells.Item(86, 4).Value2 = 0
$ws.Cells.Item(86, 5).Value2 = 0
$ws.Cells.Item(87, 4).Value2 = 0
$ws.Cells.Item(87, 5).Value2 = 0
$ws.Cells.Item(88, 4).Value2 = 0
$ws.Cells.Item(88, 5).Value2 = 0
$ws.Cells.Item(89, 4).Value2 = 126000
$ws.Cells.Item(89, 5).Value2 = 117500
$ws.Cells.Item(92, 4).Value2 = 0
$ws.Cells.Item(92, 5).Value2 = 0
$ws.Cells.Item(93, 4).Value2 = 0
$ws.Cells.Item(93, 5).Value2 = 0
$ws.Cells.Item(94, 4).Value2 = -50200
$ws.Cells.Item(94, 5).Value2 = -83300
$ws.Cells.Item(96, 4).Value2 = 0
$ws.Cells.Item(96, 5).Value2 = 0
$ws.Cells.Item(97, 4).Value2 = 0
$ws.Cells.Item(97, 5).Value2 = 0
$ws.Cells.Item(98, 4).Value2 = 0
$ws.Cells.Item(98, 5).Value2 = 0
$ws.Cells.Item(99, 4).Value2 = 0
$ws.Cells.Item(99, 5).Value2 = 0
$ws.Cells.Item(100, 4).Value2 = -60600
$ws.Cells.Item(100, 5).Value2 = -16500
$ws.Cells.Item(101, 4).Value2 = -4200
$ws.Cells.Item(101, 5).Value2 = 700
$ws.Cells.Item(102, 4).Value2 = 11000
$ws.Cells.Item(102, 5).Value2 = 18400

# Row 91 (Capital Expenditures) has restated historical values beyond the simple shift
$ws.Cells.Item(91, 4).Value2 = -42600
$ws.Cells.Item(91, 5).Value2 = -56600
$ws.Cells.Item(91, 6).Value2 = -49900
$ws.Cells.Item(91, 7).Value2 = -44200
$ws.Cells.Item(91, 8).Value2 = -39300
$ws.Cells.Item(91, 9).Value2 = -39000
$ws.Cells.Item(91, 10).Value2 = -46300
$ws.Cells.Item(91, 11).Value2 = -43200
$ws.Cells.Item(91, 12).Value2 = -44900
$ws.Cells.Item(91, 13).Value2 = -52800
